# Apply the authors'-reply changes to the gathering workbook.
#
# Summary of the change (per the commit "Add changes based on authors' replies"):
#  - "Deployment" sub-headers renamed: "Before" -> "Snapshot", "During" -> "Continuous"
#    (these headers live in row 2, columns M/N, identically on both sheets)
#  - SPADE's Granularity cell gains an extra item: "...Stack Trace" -> "...Stack Trace, Env. Var."
#  - CPL's row (row 6) switches every capability column (D:R) from "not applicable" (✗)
#    to the new "—" marker, on the "script" sheet.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("binary", "script")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Range("M2").Value = "Snapshot"
    $sheet.Range("N2").Value = "Continuous"
}

$ws = $wb.Worksheets.Item("script")

# SPADE (row 21): granularity list grows by one item.
$ws.Range("C21").Value = "Functions, Returns, Arguments, Stack Trace, Env. Var."

# CPL (row 6): replace the "✗" marks across every tracked column with "—".
$ws.Range("D6:R6").Value = "—"
